$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at 1319-1320; this pushes the existing
# rows 1319:1383 down to 1321:1385 and extends the used range to R1385.
$ws.Rows("1319:1320").Insert()

# New row 1319: Escarola / Provincia del Elquí
$ws.Cells.Item(1319, 1).Value = 5
$ws.Cells.Item(1319, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(1319, 3).Value = "Maule"
$ws.Cells.Item(1319, 4).Value = 44753
$ws.Cells.Item(1319, 5).Value = 7
$ws.Cells.Item(1319, 6).Value = 100112033
$ws.Cells.Item(1319, 7).Value = "Lechuga"
$ws.Cells.Item(1319, 8).Value = "Escarola"
$ws.Cells.Item(1319, 9).Value = "Primera"
$ws.Cells.Item(1319, 10).Value = 700
$ws.Cells.Item(1319, 11).Value = 9000
$ws.Cells.Item(1319, 12).Value = 9000
$ws.Cells.Item(1319, 13).Value = 9000
$ws.Cells.Item(1319, 14).Value = "`$/caja 15 unidades"
$ws.Cells.Item(1319, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(1319, 16).Value = 600
$ws.Cells.Item(1319, 17).Value = 15
$ws.Cells.Item(1319, 18).Value = "Hortaliza"

# New row 1320: Marina / Región del Maule
$ws.Cells.Item(1320, 1).Value = 5
$ws.Cells.Item(1320, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(1320, 3).Value = "Maule"
$ws.Cells.Item(1320, 4).Value = 44753
$ws.Cells.Item(1320, 5).Value = 7
$ws.Cells.Item(1320, 6).Value = 100112033
$ws.Cells.Item(1320, 7).Value = "Lechuga"
$ws.Cells.Item(1320, 8).Value = "Marina"
$ws.Cells.Item(1320, 9).Value = "Primera"
$ws.Cells.Item(1320, 10).Value = 500
$ws.Cells.Item(1320, 11).Value = 7000
$ws.Cells.Item(1320, 12).Value = 7000
$ws.Cells.Item(1320, 13).Value = 7000
$ws.Cells.Item(1320, 14).Value = "`$/caja 18 unidades"
$ws.Cells.Item(1320, 15).Value = "Región del Maule"
$ws.Cells.Item(1320, 16).Value = 389
$ws.Cells.Item(1320, 17).Value = 18
$ws.Cells.Item(1320, 18).Value = "Hortaliza"
